$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Target state (per diff):
#   Column D (Cotacao)      -> updated exchange-rate numbers
#   Column E (Preco Base R) -> becomes a plain NUMBER (was a shared-string
#                               text cell referencing sharedStrings 17-23)
#   Column G (Preco Final)  -> becomes a shared-string TEXT cell (was a
#                               plain number) whose text is the new value
#                               that used to live in sharedStrings 17-23
#
# To reproduce this without leaving any stray style/numberFormat residue in
# styles.xml, we:
#   1. Cut() the existing text cell from E down into G (this physically
#      relocates the already-"shared string typed" cell, so G becomes
#      t="s" and E becomes empty, with zero style side effects).
#   2. Overwrite the (now relocated) G cell's text via a throw-away
#      ="literal" formula, then Copy/PasteSpecial(values) it onto itself --
#      this "bakes" the formula into a literal shared-string value without
#      ever touching NumberFormat/Style (which is what would otherwise
#      create unused cellXfs/numFmt entries).
#   3. Write the new plain numeric values into E.
#   4. Write the new plain numeric values into D.
# ---------------------------------------------------------------------------

$newD = @{
    2 = 5.2539
    3 = 6.165912218000001
    4 = 5.2539
    5 = 5.2539
    6 = 6.165912218000001
    7 = 5.2539
    8 = 296.46
}

$newE = @{
    2 = 5253.847461
    3 = 27746.604981
    4 = 4728.457461
    5 = 4197.8661
    6 = 18497.736654
    7 = 2524.393872
    8 = 5929.2
}

$newGText = @{
    2 = "7355.39"
    3 = "55493.21"
    4 = "8038.38"
    5 = "7136.37"
    6 = "35145.70"
    7 = "5048.79"
    8 = "6818.58"
}

# Step 1: relocate the existing shared-string cell from column E to column G
# (row by row), freeing column E and turning column G into a text cell.
foreach ($r in 2..8) {
    $ws.Cells.Item($r, 5).Cut($ws.Cells.Item($r, 7))
}

# Step 2: replace each relocated G cell's text with the new value, keeping
# it a literal shared-string cell (no formula left behind, no style churn).
foreach ($r in 2..8) {
    $ws.Cells.Item($r, 7).Formula = '="' + $newGText[$r] + '"'
}
$ws.Range("G2:G8").Copy()
$ws.Range("G2:G8").PasteSpecial(-4163)  # xlPasteValues

# Step 3: write the new numeric "Preco Base Reais" values into column E.
foreach ($r in 2..8) {
    $ws.Cells.Item($r, 5).Value = $newE[$r]
}

# Step 4: write the new "Cotacao" exchange-rate values into column D.
foreach ($r in 2..8) {
    $ws.Cells.Item($r, 4).Value = $newD[$r]
}
